# Applies "more work on debugger issues" commit:
#  - TODO sheet: add a "Statis" (status) column between Task and Example,
#    and mark the "User defined functions" task as Done.
#  - Keywords sheet: document the new FOR/NEXT statement (BTOKEN_FOR,
#    VARINDEX, START/END/STEP EXPRINDEX), and remove the stray filler
#    row that had the long run of "s" characters.
#  - refresh the two sheets' selections to where the author left off.

$wb = $excel.ActiveWorkbook

$todo = $wb.Worksheets.Item("TODO")
$kw   = $wb.Worksheets.Item("Keywords")

# ---------------------------------------------------------------
# TODO sheet: insert a new column C ("Statis") before the existing
# "Example" column (which slides from C to D).
# ---------------------------------------------------------------
$todo.Columns.Item(3).Insert()

$todo.Range("C1").Value = "Statis"
$todo.Range("C2").Value = "Done"

# Match the width already used for column B.
$todo.Columns.Item(3).ColumnWidth = $todo.Columns.Item(2).ColumnWidth()

# The "Example" header that just moved into D1 gets its own bold run
# (distinct font entry from the Priority/Task/Statis bold header style).
$todo.Range("D1").Font.Bold = $true
$todo.Range("D1").Font.Size = 11

$todo.Range("C3").Select() | Out-Null

# ---------------------------------------------------------------
# Keywords sheet: new FOR / NEXT rows documenting the for-loop
# bytecode layout, inserted right after the DEF row (row 13).
# ---------------------------------------------------------------
$kw.Range("B14").Value = "FOR"
$kw.Range("C14").Value = "X"
$kw.Range("D14").Value = "X"
$kw.Range("E14").Value = "X"
$kw.Range("F14").Value = "X"
$kw.Range("G14").Value = "BTOKEN_FOR"

$kw.Range("H14").Value = "VARINDEX"
$kw.Range("L14").Value = "START: EXPRINDEX"
$kw.Range("P14").Value = "END:EXPRINDEX"
$kw.Range("T14").Value = "STEP:EXPRINDEX"

$kw.Range("H14:K14").Merge() | Out-Null
$kw.Range("L14:O14").Merge() | Out-Null
$kw.Range("P14:S14").Merge() | Out-Null
$kw.Range("T14:W14").Merge() | Out-Null

$kw.Range("H14:K14").HorizontalAlignment = -4108
$kw.Range("L14:O14").HorizontalAlignment = -4108
$kw.Range("P14:S14").HorizontalAlignment = -4108
$kw.Range("T14:W14").HorizontalAlignment = -4108

$kw.Range("B15").Value = "NEXT"

# Drop the leftover filler row ("ssssss...") that used to sit at G30.
$kw.Range("G30").ClearContents() | Out-Null

$todo.Activate() | Out-Null
$kw.Activate() | Out-Null
$kw.Range("I18").Select() | Out-Null
